# Fix for nodal real/reactive power balance set (Nm1) bug.
# Updates decision-variable values in columns B:F for rows 5-52
# (t = 4 .. 51) per the corrected simulation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("B5").Value = [double]"0.06332391581074477"
$ws.Range("C5").Value = [double]"0.06355916260465252"
$ws.Range("D5").Value = [double]"0.06382514201487238"
$ws.Range("E5").Value = [double]"0.06418211644711853"
$ws.Range("F5").Value = [double]"0.053771550193894437"
# Row 6
$ws.Range("B6").Value = [double]"0.044494579503654264"
$ws.Range("C6").Value = [double]"0.04471003908686566"
$ws.Range("D6").Value = [double]"0.04494275373159534"
$ws.Range("E6").Value = [double]"0.04529312148249381"
$ws.Range("F6").Value = [double]"0.03525584559914957"
# Row 7
$ws.Range("B7").Value = [double]"0.009997973752444302"
$ws.Range("C7").Value = [double]"0.010357918479412587"
$ws.Range("D7").Value = [double]"0.010491231652976257"
$ws.Range("E7").Value = [double]"0.010611213653251754"
$ws.Range("F7").Value = [double]"0.01042457504513188"
# Row 8
$ws.Range("B8").Value = [double]"0.008831362554646203"
$ws.Range("C8").Value = [double]"0.008491205038374274"
$ws.Range("D8").Value = [double]"0.008391156630300782"
$ws.Range("E8").Value = [double]"0.008277781311372977"
$ws.Range("F8").Value = [double]"0.008091129549612984"
# Row 9
$ws.Range("B9").Value = [double]"0.00500057137543115"
$ws.Range("C9").Value = [double]"0.005180603940082182"
$ws.Range("D9").Value = [double]"0.005247282880301897"
$ws.Range("E9").Value = [double]"0.0053072940254593505"
$ws.Range("F9").Value = [double]"0.005213943379754853"
# Row 10
$ws.Range("B10").Value = [double]"0.011369152118397572"
$ws.Range("C10").Value = [double]"0.011584617583574184"
$ws.Range("D10").Value = [double]"0.011728123964554063"
$ws.Range("E10").Value = [double]"0.011945571284034685"
$ws.Range("F10").Value = [double]"0.005758649387649733"
# Row 11
$ws.Range("B11").Value = [double]"0.018121249006797863"
$ws.Range("C11").Value = [double]"0.01758111310566857"
$ws.Range("D11").Value = [double]"0.017470218260015515"
$ws.Range("E11").Value = [double]"0.01742298958043143"
$ws.Range("F11").Value = [double]"0.013854953261631183"
# Row 12
$ws.Range("B12").Value = [double]"0.0050006681096306"
$ws.Range("C12").Value = [double]"0.005180707773371025"
$ws.Range("D12").Value = [double]"0.005247389408775227"
$ws.Range("E12").Value = [double]"0.005307403011427218"
$ws.Range("F12").Value = [double]"0.005214048529859005"
# Row 13
$ws.Range("B13").Value = [double]"0.00500063079358504"
$ws.Range("C13").Value = [double]"0.005180667720484873"
$ws.Range("D13").Value = [double]"0.005247348317387458"
$ws.Range("E13").Value = [double]"0.005307360973713176"
$ws.Range("F13").Value = [double]"0.005214007963186004"
# Row 14
$ws.Range("B14").Value = [double]"0.005000309291350155"
$ws.Range("C14").Value = [double]"0.005180322639699659"
$ws.Range("D14").Value = [double]"0.005246994289687175"
$ws.Range("E14").Value = [double]"0.005306998793222455"
$ws.Range("F14").Value = [double]"0.005213658456012855"
# Row 15
$ws.Range("B15").Value = [double]"0.012978671719861955"
$ws.Range("C15").Value = [double]"0.014664902722496813"
$ws.Range("D15").Value = [double]"0.015422012773512558"
$ws.Range("E15").Value = [double]"0.01628853837341492"
$ws.Range("F15").Value = [double]"0.015753518592662212"
# Row 16
$ws.Range("B16").Value = [double]"0.005001219967288982"
$ws.Range("C16").Value = [double]"0.005181300067762836"
$ws.Range("D16").Value = [double]"0.0052479970394988"
$ws.Range("E16").Value = [double]"0.005308024609028374"
$ws.Range("F16").Value = [double]"0.005214648510507081"
# Row 17
$ws.Range("B17").Value = [double]"0.0026953431676274807"
$ws.Range("C17").Value = [double]"0.003271018949616974"
$ws.Range("D17").Value = [double]"0.003549848770650065"
$ws.Range("E17").Value = [double]"0.003890148426135566"
$ws.Range("F17").Value = [double]"0.0036842932696121177"
# Row 18
$ws.Range("B18").Value = [double]"0.0025000759738077485"
$ws.Range("C18").Value = [double]"0.0025900908756318476"
$ws.Range("D18").Value = [double]"0.0026234299437720363"
$ws.Range("E18").Value = [double]"0.0026534352054738727"
$ws.Range("F18").Value = [double]"0.002606760378829247"
# Row 19
$ws.Range("B19").Value = [double]"0.001988491321734766"
$ws.Range("C19").Value = [double]"0.0021559703211032205"
$ws.Range("D19").Value = [double]"0.002218305446441517"
$ws.Range("E19").Value = [double]"0.002274100888020736"
$ws.Range("F19").Value = [double]"0.0021280208361198543"
# Row 20
$ws.Range("B20").Value = [double]"0.0034773596365156776"
$ws.Range("C20").Value = [double]"0.004725704210333591"
$ws.Range("D20").Value = [double]"0.005320216020576618"
$ws.Range("E20").Value = [double]"0.006040594262941408"
$ws.Range("F20").Value = [double]"0.0057964693147216535"
# Row 21
$ws.Range("B21").Value = [double]"0.002500173732131745"
$ws.Range("C21").Value = [double]"0.002590195808220757"
$ws.Range("D21").Value = [double]"0.002623537600091275"
$ws.Range("E21").Value = [double]"0.002653545345321716"
$ws.Range("F21").Value = [double]"0.002606866642118722"
# Row 22
$ws.Range("B22").Value = [double]"0.0025001977519201273"
$ws.Range("C22").Value = [double]"0.002590221589456038"
$ws.Range("D22").Value = [double]"0.002623564049706595"
$ws.Range("E22").Value = [double]"0.0026535724039902574"
$ws.Range("F22").Value = [double]"0.002606892754028105"
# Row 23
$ws.Range("B23").Value = [double]"0.0024998716494942816"
$ws.Range("C23").Value = [double]"0.002589871571107919"
$ws.Range("D23").Value = [double]"0.002623204956427134"
$ws.Range("E23").Value = [double]"0.0026532050412667596"
$ws.Range("F23").Value = [double]"0.0026065382459574075"
# Row 24
$ws.Range("B24").Value = [double]"0.002024662168522598"
$ws.Range("C24").Value = [double]"0.0020868091920567267"
$ws.Range("D24").Value = [double]"0.0021279751529137613"
$ws.Range("E24").Value = [double]"0.002183661074367119"
$ws.Range("F24").Value = [double]"0.0014055498550281777"
# Row 25
$ws.Range("B25").Value = [double]"0.00011779779462322102"
$ws.Range("C25").Value = [double]"0.00012643260026237716"
$ws.Range("D25").Value = [double]"0.00012970818599585172"
$ws.Range("E25").Value = [double]"0.00013269203548402987"
$ws.Range("F25").Value = [double]"0.00012806515596165523"
# Row 26
$ws.Range("B26").Value = [double]"8.01271931734294e-5"
$ws.Range("C26").Value = [double]"7.783592052828577e-5"
$ws.Range("D26").Value = [double]"7.804415762771203e-5"
$ws.Range("E26").Value = [double]"7.865648042758507e-5"
$ws.Range("F26").Value = [double]"7.450178060733178e-5"
# Row 27
$ws.Range("B27").Value = [double]"2.9473771198790548e-5"
$ws.Range("C27").Value = [double]"3.163491690433745e-5"
$ws.Range("D27").Value = [double]"3.245483642782512e-5"
$ws.Range("E27").Value = [double]"3.32018397263991e-5"
$ws.Range("F27").Value = [double]"3.204208523096305e-5"
# Row 28
$ws.Range("B28").Value = [double]"0.00012554357077840337"
$ws.Range("C28").Value = [double]"0.00013090660220814636"
$ws.Range("D28").Value = [double]"0.00013432468139175897"
$ws.Range("E28").Value = [double]"0.00013938119969232358"
$ws.Range("F28").Value = [double]"3.554053815717879e-5"
# Row 29
$ws.Range("B29").Value = [double]"0.0003209791638092843"
$ws.Range("C29").Value = [double]"0.0003124925463587206"
$ws.Range("D29").Value = [double]"0.00031447082366101506"
$ws.Range("E29").Value = [double]"0.0003206462652343732"
$ws.Range("F29").Value = [double]"0.0002126923455470125"
# Row 30
$ws.Range("B30").Value = [double]"2.9485639292648346e-5"
$ws.Range("C30").Value = [double]"3.164808646761617e-5"
$ws.Range("D30").Value = [double]"3.2468630426003706e-5"
$ws.Range("E30").Value = [double]"3.3216343344529106e-5"
$ws.Range("F30").Value = [double]"3.205404154674865e-5"
# Row 31
$ws.Range("B31").Value = [double]"2.948540053924244e-5"
$ws.Range("C31").Value = [double]"3.164782096790132e-5"
$ws.Range("D31").Value = [double]"3.2468354531195046e-5"
$ws.Range("E31").Value = [double]"3.321605786101328e-5"
$ws.Range("F31").Value = [double]"3.205377093348442e-5"
# Row 32
$ws.Range("B32").Value = [double]"2.9485400539242442e-5"
$ws.Range("C32").Value = [double]"3.1647820967901323e-5"
$ws.Range("D32").Value = [double]"3.246835453119505e-5"
$ws.Range("E32").Value = [double]"3.321605786101327e-5"
$ws.Range("F32").Value = [double]"3.205377093348443e-5"
# Row 34
$ws.Range("B34").Value = [double]"1.060471486463274"
$ws.Range("C34").Value = [double]"1.0604493771525947"
$ws.Range("D34").Value = [double]"1.0604386305180635"
$ws.Range("E34").Value = [double]"1.0604258243415199"
$ws.Range("F34").Value = [double]"1.0604929421723979"
# Row 35
$ws.Range("B35").Value = [double]"1.0607372295225344"
$ws.Range("C35").Value = [double]"1.060731369416348"
$ws.Range("D35").Value = [double]"1.0607291990001404"
$ws.Range("E35").Value = [double]"1.0607272456225454"
$ws.Range("F35").Value = [double]"1.0607302842086839"
# Row 36
$ws.Range("B36").Value = [double]"1.0603071426960828"
$ws.Range("C36").Value = [double]"1.0602791167444305"
$ws.Range("D36").Value = [double]"1.060266178756688"
$ws.Range("E36").Value = [double]"1.0602514003600656"
$ws.Range("F36").Value = [double]"1.0603215860883013"
# Row 37
$ws.Range("B37").Value = [double]"1.060378363542893"
$ws.Range("C37").Value = [double]"1.0603528875405295"
$ws.Range("D37").Value = [double]"1.0603405059717077"
$ws.Range("E37").Value = [double]"1.0603256995976507"
$ws.Range("F37").Value = [double]"1.0604319388490702"
# Row 38
$ws.Range("B38").Value = [double]"1.0600940315335743"
$ws.Range("C38").Value = [double]"1.0600592674234022"
$ws.Range("D38").Value = [double]"1.0600399429619058"
$ws.Range("E38").Value = [double]"1.0600152237791725"
$ws.Range("F38").Value = [double]"1.0601488649084456"
# Row 39
$ws.Range("B39").Value = [double]"1.059880364519082"
$ws.Range("C39").Value = [double]"1.0598379079476403"
$ws.Range("D39").Value = [double]"1.0598157344184105"
$ws.Range("E39").Value = [double]"1.0597884510709572"
$ws.Range("F39").Value = [double]"1.0599260809005766"
# Row 40
$ws.Range("B40").Value = [double]"1.0599296837797423"
$ws.Range("C40").Value = [double]"1.05988900273612"
$ws.Range("D40").Value = [double]"1.0598674868104008"
$ws.Range("E40").Value = [double]"1.0598407953063955"
$ws.Range("F40").Value = [double]"1.0599775044906765"
# Row 41
$ws.Range("B41").Value = [double]"1.0598889467340868"
$ws.Range("C41").Value = [double]"1.0598467991336291"
$ws.Range("D41").Value = [double]"1.0598247400382979"
$ws.Range("E41").Value = [double]"1.0597975596814475"
$ws.Range("F41").Value = [double]"1.0599350293034706"
# Row 42
$ws.Range("B42").Value = [double]"1.0607232105806008"
$ws.Range("C42").Value = [double]"1.060719489586395"
$ws.Range("D42").Value = [double]"1.0607166981199432"
$ws.Range("E42").Value = [double]"1.0607131578303646"
$ws.Range("F42").Value = [double]"1.0607192014172477"
# Row 43
$ws.Range("B43").Value = [double]"0.005036026207556386"
$ws.Range("C43").Value = [double]"0.004132669373122418"
$ws.Range("D43").Value = [double]"0.003666872680988475"
$ws.Range("E43").Value = [double]"0.0030622165506399724"
$ws.Range("F43").Value = [double]"0.003089824594449595"
# Row 44
$ws.Range("B44").Value = [double]"0.0023055277210955417"
$ws.Range("C44").Value = [double]"0.0019098431845174946"
$ws.Range("D44").Value = [double]"0.0016976849770697743"
$ws.Range("E44").Value = [double]"0.0014173930663445257"
$ws.Range("F44").Value = [double]"0.0015298780289751922"
# Row 45
$ws.Range("B45").Value = [double]"0.003012533444693624"
$ws.Range("C45").Value = [double]"0.0030250965421749765"
$ws.Range("D45").Value = [double]"0.003029452510812138"
$ws.Range("E45").Value = [double]"0.00303369127449974"
$ws.Range("F45").Value = [double]"0.003085792495030601"
# Row 46
$ws.Range("B46").Value = [double]"0.0014900733847792735"
$ws.Range("C46").Value = [double]"0.0015050655296811194"
$ws.Range("D46").Value = [double]"0.0015097221534455907"
$ws.Range("E46").Value = [double]"0.0015140737374069323"
$ws.Range("F46").Value = [double]"0.00154293987880631"
# Row 47
$ws.Range("B47").Value = [double]"0.0013713150393736503"
$ws.Range("C47").Value = [double]"0.0012268543288419804"
$ws.Range("D47").Value = [double]"0.0012370504002260617"
$ws.Range("E47").Value = [double]"0.0013345124680544814"
$ws.Range("F47").Value = [double]"-8.918345262159885e-9"
# Row 48
$ws.Range("B48").Value = [double]"0.0004524429015939975"
$ws.Range("C48").Value = [double]"0.000592044405570547"
$ws.Range("D48").Value = [double]"0.0006810200770422861"
$ws.Range("E48").Value = [double]"0.0008602912905415712"
$ws.Range("F48").Value = [double]"-8.918550613658096e-9"
# Row 49
$ws.Range("B49").Value = [double]"4.159894950410076e-8"
$ws.Range("C49").Value = [double]"4.1599791410278856e-8"
$ws.Range("D49").Value = [double]"4.160017238244907e-8"
$ws.Range("E49").Value = [double]"4.16005326142097e-8"
$ws.Range("F49").Value = [double]"0.004665508895185284"
# Row 50
$ws.Range("B50").Value = [double]"4.159337281592296e-8"
$ws.Range("C50").Value = [double]"4.1589594244711796e-8"
$ws.Range("D50").Value = [double]"4.158872287967354e-8"
$ws.Range("E50").Value = [double]"4.158788425774589e-8"
$ws.Range("F50").Value = [double]"0.0023335088953943058"
# Row 51
$ws.Range("B51").Value = [double]"0.012966455499037069"
$ws.Range("C51").Value = [double]"0.014131923322182834"
$ws.Range("D51").Value = [double]"0.015307077412742452"
$ws.Range("E51").Value = [double]"0.016574820467359877"
# Row 52
$ws.Range("B52").Value = [double]"0.006263401974016596"
$ws.Range("C52").Value = [double]"0.006825800380788358"
$ws.Range("D52").Value = [double]"0.0074727256763754985"
$ws.Range("E52").Value = [double]"0.00828995862566972"
